$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New schools data: School Name, Achievement(C), Growth(D), Reading(E), Math(F), Science(G),
# SchoolSize(H), TeachRatio(I), White(J), Black(K), Hispanic(L), Asian(M), Mixed/other(N)
# Performance Score (B) is computed with the same formula used for the existing rows.
$newRows = @(
    @{ Row = 21; Name = "Elon Park";              C = 88.5; D = 87;   E = 83.9; F = 92;   G = 93.4; H = 1004; I = 17.6; J = 49.6; K = 7.8;  L = 5.7;  M = 34.7; N = 2.2 },
    @{ Row = 22; Name = "Hawk Ridge";              C = 85.9; D = 86.7; E = 82.2; F = 86.3; G = 95;   H = 881;  I = 17.2; J = 59.1; K = 9.3;  L = 9.3;  M = 19.3; N = 3   },
    @{ Row = 23; Name = "Weddington Elementary";   C = 95.1; D = 89.8; E = 93;   F = 95;   G = 95;   H = 851;  I = 17.4; J = 73.5; K = 15.4; L = 3.9;  M = 5.4;  N = 1.8 },
    @{ Row = 24; Name = "Olde Providence";         C = 81.8; D = 72.2; E = 79.4; F = 85.3; G = 78;   H = 715;  I = 19.4; J = 73.8; K = 8.9;  L = 5.3;  M = 9.6;  N = 2.4 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
}

# Performance score = (0.8 * Achievement Score) + (0.2 * Growth Score), same formula used for the
# rest of the data (B3:B20). Assigning it across the whole new block at once lets it fill down as a
# single formula (relative references adjust automatically per row, same as the existing column).
$ws.Range("B21:B24").Formula = "=(0.8*C21)+(0.2*D21)"

# Match the centered-number style (s="2") used by the rest of the data rows
$ws.Range("B21:N24").HorizontalAlignment = -4108

# Update the active selection to match the new extent of the data
$ws.Range("O19:O26").Select() | Out-Null
